$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing contents of row 1 and 2 (B1:C2) that are being removed
$ws.Range("B1:C2").Clear()

# Update header and value
$ws.Range("A1").Value = "zip_code"
$ws.Range("A2").Value = 94553

# Update selection to match target state
$ws.Range("A2").Select()
